# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型" sheets
# to reflect refreshed scrape data, per the commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 45
$ws1.Range("F5").Value = 19
$ws1.Range("F7").Value = 14652
$ws1.Range("F9").Value = 668
$ws1.Range("F10").Value = 15158
$ws1.Range("F12").Value = 8597
$ws1.Range("F13").Value = 314
$ws1.Range("F15").Value = 60
$ws1.Range("F16").Value = 176
$ws1.Range("F20").Value = 7
$ws1.Range("F21").Value = 20
$ws1.Range("F24").Value = 1070
$ws1.Range("F26").Value = 8
$ws1.Range("F31").Value = 21
$ws1.Range("F32").Value = 24
$ws1.Range("F33").Value = 227
$ws1.Range("F37").Value = 5307

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 45
$ws4.Range("F5").Value = 19
$ws4.Range("F7").Value = 14652
$ws4.Range("F9").Value = 668
$ws4.Range("F10").Value = 15158
$ws4.Range("F12").Value = 8597
$ws4.Range("F13").Value = 314
$ws4.Range("F16").Value = 60
$ws4.Range("F17").Value = 176
$ws4.Range("F21").Value = 7
$ws4.Range("F22").Value = 20
$ws4.Range("F25").Value = 1070
$ws4.Range("F27").Value = 8
$ws4.Range("F34").Value = 21
$ws4.Range("F35").Value = 24
$ws4.Range("F36").Value = 227
$ws4.Range("F40").Value = 5307

$wb.Save()
